$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D hold numeric-looking values but must remain plain text,
# matching the original inline-string cell type. Force text format first
# so Excel does not auto-convert them to numbers.

$ws.Range('D2').Value = '69.324.01'
$ws.Range('E2').Value = '  -1.18%  '
$ws.Range('D3').Value = '2.516.56'
$ws.Range('E3').Value = '  -0.35%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '572.42'
$ws.Range('E5').Value = '  -0.42%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '166.62'
$ws.Range('E6').Value = '  -1.98%  '
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.517'
$ws.Range('E8').Value = '  +1.58%  '
$ws.Range('D9').Value = '2.515.22'
$ws.Range('E9').Value = '  -0.39%  '
$ws.Range('E10').Value = '  +0.01%  '
$ws.Range('E11').Value = '  -0.48%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.355'
$ws.Range('E12').Value = '  +3.50%  '
$ws.Range('E13').Value = '  +2.50%  '
$ws.Range('D14').Value = '2.980.18'
$ws.Range('E14').Value = '  -0.33%  '
$ws.Range('D15').Value = '69.247.25'
$ws.Range('E15').Value = '  -1.08%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0000176'
$ws.Range('E16').Value = '  -2.09%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '24.87'
$ws.Range('E17').Value = '  +0.11%  '
$ws.Range('D18').Value = '2.518.63'
$ws.Range('E18').Value = '  +0.10%  '
$ws.Range('E19').Value = '  -0.64%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.62'
$ws.Range('E20').Value = '  +1.14%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '348.44'
$ws.Range('E21').Value = '  -1.47%  '
$ws.Range('E22').Value = '  -0.12%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '2.00'
$ws.Range('E23').Value = '  +1.14%  '
$ws.Range('E24').Value = '  +0.06%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '70.27'
$ws.Range('E25').Value = '  +2.11%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.97'
$ws.Range('E26').Value = '  -2.07%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.92'
$ws.Range('E27').Value = '  -3.25%  '
$ws.Range('E28').Value = '  -0.33%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.995'
$ws.Range('E29').Value = '  -1.02%  '
$ws.Range('D30').Value = '0.0₃0894'
$ws.Range('E30').Value = '  -1.84%  '
$ws.Range('E31').Value = '  +0.07%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '462.60'
$ws.Range('E32').Value = '  -3.51%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.24'
$ws.Range('E33').Value = '  -3.98%  '
$ws.Range('E34').Value = '  -1.30%  '
$ws.Range('E35').Value = '  +0.06%  '
$ws.Range('E36').Value = '  +1.29%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '157.64'
$ws.Range('E37').Value = '  +0.22%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '19.01'
$ws.Range('E38').Value = '  +0.92%  '
$ws.Range('E39').Value = '  -0.07%  '
$ws.Range('E40').Value = '  +0.00%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '4.73'
$ws.Range('E41').Value = '  +0.45%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.319'
$ws.Range('E42').Value = '  -0.13%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.61'
$ws.Range('E43').Value = '  -2.51%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '38.12'
$ws.Range('E44').Value = '  -0.50%  '
$ws.Range('B45').Value = 'dogwifhat'
$ws.Range('C45').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.26'
$ws.Range('E45').Value = '  -5.06%  '
$ws.Range('B46').Value = 'ImmutableX'
$ws.Range('C46').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.13'
$ws.Range('E46').Value = '  -13.31%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '141.66'
$ws.Range('E47').Value = '  -0.25%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.524'
$ws.Range('E48').Value = '  +0.20%  '
$ws.Range('E49').Value = '  -1.45%  '
$ws.Range('E50').Value = '  -0.19%  '
$ws.Range('B51').Value = 'Mantle'
$ws.Range('C51').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.579'
$ws.Range('E51').Value = '  -3.21%  '
